$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.048.59'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.65%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.357.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +6.38%  '

$ws.Range("E4").Value = '  -0.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.62%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.637'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.44%  '

$ws.Range("E8").Value = '  -0.29%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.631'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.17'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.81%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0935'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.88'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.04'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +9.88%  '

$ws.Range("E14").Value = '  +1.09%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.43'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +9.18%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.714.79'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.464.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +10.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.025.47'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.70%  '

$ws.Range("E19").Value = '  +2.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '75.18'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.97%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.39%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '252.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +10.16%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.86'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.60%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.32%  '

$ws.Range("E27").Value = '  +0.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.23'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.49%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '38.51'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.40%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.58'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.36%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.20'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.69%  '

$ws.Range("E32").Value = '  -2.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0915'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.77%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.90%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.131'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.96'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.25%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0375'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.49%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.04'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.54%  '

$ws.Range("E39").Value = '  +0.43%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.73'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.57%  '

$ws.Range("E41").Value = '  +14.38%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '71.60'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.229'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.27%  '

$ws.Range("E44").Value = '  -0.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.25'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.64%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.60'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.92%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '110.05'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.52%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.28'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.41%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0996'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.493.98'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.98%  '
